$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text representation (some values
# look numeric, e.g. "28.60", and would otherwise be coerced to numbers,
# dropping formatting such as trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.738.32"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").Value = "2.654.18"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "598.06"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("D6").Value = "157.25"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "0.635"
$ws.Range("E8").Value = "  +1.93%  "

$ws.Range("E9").Value = "  -2.93%  "

$ws.Range("E10").Value = "  -1.36%  "

$ws.Range("D11").Value = "5.81"
$ws.Range("E11").Value = "  -0.75%  "

$ws.Range("D12").Value = "0.157"
$ws.Range("E12").Value = "  +1.21%  "

$ws.Range("D13").Value = "28.60"
$ws.Range("E13").Value = "  -3.36%  "

$ws.Range("E14").Value = "  -2.28%  "

$ws.Range("D15").Value = "3.131.52"
$ws.Range("E15").Value = "  -0.73%  "

$ws.Range("D16").Value = "65.585.97"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").Value = "2.602.23"
$ws.Range("E17").Value = "  -2.66%  "

$ws.Range("E18").Value = "  -1.97%  "

$ws.Range("D19").Value = "4.75"
$ws.Range("E19").Value = "  -1.58%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").Value = "69.36"
$ws.Range("E23").Value = "  -0.51%  "

$ws.Range("D24").Value = "1.79"
$ws.Range("E24").Value = "  +9.40%  "

$ws.Range("D25").Value = "0.0000111"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("D26").Value = "9.59"
$ws.Range("E26").Value = "  -2.52%  "

$ws.Range("D27").Value = "1.61"
$ws.Range("E27").Value = "  +0.94%  "

$ws.Range("D28").Value = "568.30"
$ws.Range("E28").Value = "  +6.68%  "

$ws.Range("D29").Value = "0.164"
$ws.Range("E29").Value = "  -2.66%  "

$ws.Range("E32").Value = "  -0.91%  "

$ws.Range("D33").Value = "1.78"
$ws.Range("E33").Value = "  +0.86%  "

$ws.Range("D34").Value = "6.53"
$ws.Range("E34").Value = "  -0.81%  "

$ws.Range("E35").Value = "  -1.79%  "

$ws.Range("E36").Value = "  -1.47%  "

$ws.Range("D37").Value = "20.40"
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("D40").Value = "154.38"
$ws.Range("E40").Value = "  -3.26%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").Value = "161.70"
$ws.Range("E42").Value = "  -2.23%  "

$ws.Range("D43").Value = "4.08"
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("E44").Value = "  -1.15%  "

$ws.Range("E45").Value = "  -1.48%  "

$ws.Range("D46").Value = "22.57"
$ws.Range("E46").Value = "  -2.40%  "

$ws.Range("D47").Value = "0.638"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("E48").Value = "  -1.82%  "

$ws.Range("D49").Value = "0.101"
$ws.Range("E49").Value = "  -1.18%  "

$ws.Range("D50").Value = "19.72"
$ws.Range("E50").Value = "  -3.32%  "

$ws.Range("E51").Value = "  +5.64%  "

# Row 20 and 21 swapped (BitcoinCash moved up, Uniswap moved down)
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "349.73"
$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "7.44"
$ws.Range("E21").Value = "  -3.26%  "

# Row 30 and 31 swapped (Aptos moved up, Binance-PegBSC-USD moved down)
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "8.01"
$ws.Range("E30").Value = "  -1.31%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "0.995"
$ws.Range("E31").Value = "  -0.49%  "
